$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New defined name used by the antenna-calculation formulas below.
$wb.Names.Add("Fc", '=Sheet1!$B$25')

# New "Antenna calculations" section (rows 24-27).
$ws.Range("A24").Value = "Antenna calculations"
$ws.Range("A24").Style = $ws.Range("A3").Style

$ws.Range("A25").Value = "Fcarrier"
$ws.Range("B25").Value = 315
$ws.Range("C25").Value = "MHz"

$ws.Range("A26").Value = "Wave length"
$ws.Range("B26").Formula = "=3*10^8/(Fc*10^6)"
$ws.Range("B26").NumberFormat = "0.00"
$ws.Range("C26").Value = "m"

$ws.Range("A27").Value = "WL/4"
$ws.Range("B27").Formula = "=B26/4"
$ws.Range("B27").NumberFormat = "0.00"
$ws.Range("C27").Value = "m"

$ws.Range("D29").Select()
